$wb = $excel.ActiveWorkbook

# The same six cells need updating on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8468
    $ws.Range("F3").Value = 8160
    $ws.Range("F4").Value = 141
    $ws.Range("F13").Value = 196
    $ws.Range("F14").Value = 4233
    $ws.Range("F20").Value = 122
}
